# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# sheets, which hold identical duplicated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1166
    4  = 1595
    5  = 611
    6  = 1093
    7  = 12
    8  = 11382
    9  = 20
    11 = 446
    14 = 786
    15 = 12340
    16 = 13011
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
